# Adapt column header formatting to respective input file names (old/new -> FV2210/FV2304),
# freeze the header row, and wrap the data range in an Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) from "<Name>_old" / "<Name>_new" to
#        "<Name>_FV2210" / "<Name>_FV2304" respectively. Column "diff" (K1) stays.
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $baseNames[$i] + "_FV2210"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2304"
}

# --- 2. Freeze the header row (row 1) ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the data range A1:U65 in an Excel Table ("Table1") ---
$dataRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1, $null)
$table.Name = "Table1"

[void]$ws.Range("A1").Select()
